# Convert the "git://github.com/bam593/bmProjects.git" hyperlink into
# plain text preceded by "git clone ", matching the author's edit that
# turns:
#   <w:hyperlink r:id="rIdX" w:history="1">
#     <w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr>
#       <w:t>git://github.com/bam593/bmProjects.git</w:t></w:r>
#   </w:hyperlink>
# into:
#   <w:proofErr w:type="spellStart"/>
#   <w:r><w:t>git</w:t></w:r>
#   <w:proofErr w:type="spellEnd"/>
#   <w:r><w:t xml:space="preserve"> clone </w:t></w:r>
#   <w:r><w:t>git://github.com/bam593/bmProjects.git</w:t></w:r>

$d = $word.ActiveDocument

$targetAddr = "git://github.com/bam593/bmProjects.git"

# Locate the hyperlink object by its target address (robust to index shifts).
$hyperlinkObj = $null
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $candidate = $d.Hyperlinks.Item($i)
    if ($candidate.Address -eq $targetAddr) {
        $hyperlinkObj = $candidate
        break
    }
}

if ($hyperlinkObj -eq $null) {
    Write-Host "ERROR: target hyperlink not found"
} else {
    # Find which paragraph contains the hyperlink (iterate rather than using
    # Range.Paragraphs, which mis-anchors on non-whole-document ranges here).
    $hlStart = $hyperlinkObj.Range.Start
    $paraIndex = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $cand = $d.Paragraphs.Item($i)
        if ($cand.Range.Start -le $hlStart -and $hlStart -lt $cand.Range.End) {
            $paraIndex = $i
            break
        }
    }

    # Removing the hyperlink turns its run back into plain text (still
    # carrying the Hyperlink character style), but keeps it a single run.
    $hyperlinkObj.Delete()

    # Re-fetch the (now plain-text) paragraph and overwrite its whole
    # contents (including the paragraph mark) with the three target runs
    # plus the spell-check markers around "git", via a WordprocessingML
    # fragment so no Hyperlink run style survives.
    $para = $d.Paragraphs.Item($paraIndex)
    $paraRange = $para.Range

    $xmlFragment = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> clone </w:t></w:r><w:r><w:t>$targetAddr</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

    $null = $paraRange.InsertXML($xmlFragment)
}
